$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 11.921572
$ws.Range("H2").Value = 35.764716
$ws.Range("I2").Value = 0.006971694289596158
$ws.Range("J2").Value = 0.006971694289596159
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.879565666666667
$ws.Range("N2").Value = 23.638697
$ws.Range("O2").Value = 0.9977172793687663
$ws.Range("P2").Value = 0.9977172793687664
$ws.Range("Q2").Value = 93.93680942389467
$ws.Range("R2").Value = 845.431284815052
$ws.Range("S2").Value = 0.006955779859206643
$ws.Range("T2").Value = 0.006955779859206645

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 11.921572
$ws.Range("H3").Value = 35.764716
$ws.Range("I3").Value = 0.006971694289596158
$ws.Range("J3").Value = 0.006971694289596159
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.018028
$ws.Range("N3").Value = 0.054084
$ws.Range("O3").Value = 0.002282720631233623
$ws.Range("P3").Value = 0.002282720631233623
$ws.Range("Q3").Value = 0.214922100016
$ws.Range("R3").Value = 1.934298900144
$ws.Range("S3").Value = 0.00001591443038951479
$ws.Range("T3").Value = 0.00001591443038951479

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1516.953124666667
$ws.Range("H4").Value = 4550.859374
$ws.Range("I4").Value = 0.8871089682487887
$ws.Range("J4").Value = 0.8871089682487888
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.879565666666667
$ws.Range("N4").Value = 23.638697
$ws.Range("O4").Value = 0.9977172793687663
$ws.Range("P4").Value = 0.9977172793687664
$ws.Range("Q4").Value = 11952.93175906619
$ws.Range("R4").Value = 107576.3858315957
$ws.Range("S4").Value = 0.8850839463048148
$ws.Range("T4").Value = 0.885083946304815

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1516.953124666667
$ws.Range("H5").Value = 4550.859374
$ws.Range("I5").Value = 0.8871089682487887
$ws.Range("J5").Value = 0.8871089682487888
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.018028
$ws.Range("N5").Value = 0.054084
$ws.Range("O5").Value = 0.002282720631233623
$ws.Range("P5").Value = 0.002282720631233623
$ws.Range("Q5").Value = 27.34763093149066
$ws.Range("R5").Value = 246.128678383416
$ws.Range("S5").Value = 0.002025021943973883
$ws.Range("T5").Value = 0.002025021943973884

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 181.1216836666667
$ws.Range("H6").Value = 543.365051
$ws.Range("I6").Value = 0.1059193374616151
$ws.Range("J6").Value = 0.1059193374616151
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.879565666666667
$ws.Range("N6").Value = 23.638697
$ws.Range("O6").Value = 0.9977172793687663
$ws.Range("P6").Value = 0.9977172793687664
$ws.Range("Q6").Value = 1427.160200108727
$ws.Range("R6").Value = 12844.44180097855
$ws.Range("S6").Value = 0.1056775532047449
$ws.Range("T6").Value = 0.1056775532047449

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 181.1216836666667
$ws.Range("H7").Value = 543.365051
$ws.Range("I7").Value = 0.1059193374616151
$ws.Range("J7").Value = 0.1059193374616151
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.018028
$ws.Range("N7").Value = 0.054084
$ws.Range("O7").Value = 0.002282720631233623
$ws.Range("P7").Value = 0.002282720631233623
$ws.Range("Q7").Value = 3.265261713142666
$ws.Range("R7").Value = 29.387355418284
$ws.Range("S7").Value = 0.0002417842568702252
$ws.Range("T7").Value = 0.0002417842568702252
